$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 472.08694  # H92
$ws.Cells.Item(92, 9).Value = 380.4  # I92
$ws.Cells.Item(92, 11).Value = 380.4  # K92
$ws.Cells.Item(92, 13).Value = 867.6  # M92
$ws.Cells.Item(100, 8).Value = 2853.6365  # H100
$ws.Cells.Item(100, 10).Value = 3350  # J100
$ws.Cells.Item(100, 12).Value = 3350  # L100
$ws.Cells.Item(100, 14).Value = -4432  # N100
$ws.Cells.Item(101, 8).Value = 1334.1666  # H101
$ws.Cells.Item(101, 9).Value = 501.25  # I101
$ws.Cells.Item(101, 10).Value = 3000  # J101
$ws.Cells.Item(101, 11).Value = 1503.75  # K101
$ws.Cells.Item(101, 12).Value = 9000  # L101
$ws.Cells.Item(101, 13).Value = 118.25  # M101
$ws.Cells.Item(101, 14).Value = -12244  # N101
$ws.Cells.Item(116, 8).Value = 4453.1333  # H116
$ws.Cells.Item(116, 9).Value = 2496.25  # I116
$ws.Cells.Item(116, 10).Value = 5164.727  # J116
$ws.Cells.Item(116, 11).Value = 2496.25  # K116
$ws.Cells.Item(116, 12).Value = 5164.727  # L116
$ws.Cells.Item(116, 13).Value = 945.75  # M116
$ws.Cells.Item(116, 14).Value = -12048.727  # N116
$ws.Cells.Item(121, 8).Value = 1851.5  # H121
$ws.Cells.Item(121, 10).Value = 1968.3334  # J121
$ws.Cells.Item(121, 12).Value = 5905.0002  # L121
$ws.Cells.Item(121, 14).Value = -9399.0002  # N121
$ws.Cells.Item(125, 8).Value = 1272  # H125
$ws.Cells.Item(125, 9).Value = 298.66666  # I125
$ws.Cells.Item(125, 10).Value = 1856  # J125
$ws.Cells.Item(125, 11).Value = 2687.99994  # K125
$ws.Cells.Item(125, 12).Value = 16704  # L125
$ws.Cells.Item(125, 13).Value = -227.9999399999997  # M125
$ws.Cells.Item(125, 14).Value = -21624  # N125
$ws.Cells.Item(129, 8).Value = 843.3333  # H129
$ws.Cells.Item(129, 10).Value = 848.8333  # J129
$ws.Cells.Item(129, 12).Value = 2546.4999  # L129
$ws.Cells.Item(129, 14).Value = -12546.4999  # N129
$ws.Cells.Item(141, 8).Value = 4198.75  # H141
$ws.Cells.Item(141, 10).Value = 4725  # J141
$ws.Cells.Item(141, 12).Value = 14175  # L141
$ws.Cells.Item(141, 14).Value = -24535  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 2297.2778  # H102
$ws.Cells.Item(102, 9).Value = 1091.4  # I102
$ws.Cells.Item(102, 11).Value = 1091.4  # K102
$ws.Cells.Item(102, 13).Value = 530.5999999999999  # M102

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 4338.5264  # H94
$ws.Cells.Item(94, 9).Value = 1942.2  # I94
$ws.Cells.Item(94, 10).Value = 7001.1113  # J94
$ws.Cells.Item(94, 11).Value = 1942.2  # K94
$ws.Cells.Item(94, 12).Value = 7001.1113  # L94
$ws.Cells.Item(94, 13).Value = -1491.2  # M94
$ws.Cells.Item(94, 14).Value = -7903.1113  # N94

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 950  # H22
$ws.Cells.Item(22, 10).Value = 0  # J22
$ws.Cells.Item(22, 12).Value = 0  # L22
$ws.Cells.Item(22, 14).ClearContents()  # N22: was -1700
$ws.Cells.Item(31, 8).Value = 10139.863  # H31
$ws.Cells.Item(31, 9).Value = 11820.143  # I31
$ws.Cells.Item(31, 11).Value = 11820.143  # K31
$ws.Cells.Item(31, 13).Value = -11525.143  # M31
$ws.Cells.Item(34, 8).Value = 10139.863  # H34
$ws.Cells.Item(34, 9).Value = 11820.143  # I34
$ws.Cells.Item(34, 11).Value = 11820.143  # K34
$ws.Cells.Item(34, 13).Value = -11618.143  # M34
$ws.Cells.Item(105, 8).Value = 31251116  # H105
$ws.Cells.Item(105, 9).Value = 62500224  # I105
$ws.Cells.Item(105, 10).Value = 2005.5  # J105
$ws.Cells.Item(105, 11).Value = 62500224  # K105
$ws.Cells.Item(105, 12).Value = 2005.5  # L105
$ws.Cells.Item(105, 13).Value = -62498477  # M105
$ws.Cells.Item(105, 14).Value = -5499.5  # N105
$ws.Cells.Item(122, 8).Value = 1176.9231  # H122
$ws.Cells.Item(122, 9).Value = 1132.8334  # I122
$ws.Cells.Item(122, 10).Value = 1214.7142  # J122
$ws.Cells.Item(122, 11).Value = 3398.5002  # K122
$ws.Cells.Item(122, 12).Value = 3644.1426  # L122
$ws.Cells.Item(122, 13).Value = -948.5001999999999  # M122
$ws.Cells.Item(122, 14).Value = -8544.142599999999  # N122
$ws.Cells.Item(132, 8).Value = 26649.809  # H132
$ws.Cells.Item(132, 9).Value = 42935.418  # I132
$ws.Cells.Item(132, 10).Value = 4935.6665  # J132
$ws.Cells.Item(132, 11).Value = 128806.254  # K132
$ws.Cells.Item(132, 12).Value = 14806.9995  # L132
$ws.Cells.Item(132, 13).Value = -126276.254  # M132
$ws.Cells.Item(132, 14).Value = -19866.9995  # N132
$ws.Cells.Item(134, 8).Value = 1389.8334  # H134
$ws.Cells.Item(134, 9).Value = 1109.1428  # I134
$ws.Cells.Item(134, 10).Value = 1782.8  # J134
$ws.Cells.Item(134, 11).Value = 3327.4284  # K134
$ws.Cells.Item(134, 12).Value = 5348.4  # L134
$ws.Cells.Item(134, 13).Value = -792.4284000000002  # M134
$ws.Cells.Item(134, 14).Value = -10418.4  # N134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 117075.164  # H131
$ws.Cells.Item(131, 9).Value = 732  # I131
$ws.Cells.Item(131, 10).Value = 124256.836  # J131
$ws.Cells.Item(131, 11).Value = 2196  # K131
$ws.Cells.Item(131, 12).Value = 372770.508  # L131
$ws.Cells.Item(131, 13).Value = 2844  # M131
$ws.Cells.Item(131, 14).Value = -382850.508  # N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4405.7144  # H80
$ws.Cells.Item(80, 9).Value = 3775  # I80
$ws.Cells.Item(80, 10).Value = 4658  # J80
$ws.Cells.Item(80, 11).Value = 3775  # K80
$ws.Cells.Item(80, 12).Value = 4658  # L80
$ws.Cells.Item(80, 13).Value = -2777  # M80
$ws.Cells.Item(80, 14).Value = -6654  # N80
$ws.Cells.Item(83, 8).Value = 4405.7144  # H83
$ws.Cells.Item(83, 9).Value = 3775  # I83
$ws.Cells.Item(83, 10).Value = 4658  # J83
$ws.Cells.Item(83, 11).Value = 18875  # K83
$ws.Cells.Item(83, 12).Value = 23290  # L83
$ws.Cells.Item(83, 13).Value = -13883  # M83
$ws.Cells.Item(83, 14).Value = -33274  # N83
$ws.Cells.Item(97, 8).Value = 2142.8333  # H97
$ws.Cells.Item(97, 9).Value = 1227.4286  # I97
$ws.Cells.Item(97, 10).Value = 5346.75  # J97
$ws.Cells.Item(97, 11).Value = 1227.4286  # K97
$ws.Cells.Item(97, 12).Value = 5346.75  # L97
$ws.Cells.Item(97, 13).Value = -731.4286  # M97
$ws.Cells.Item(97, 14).Value = -6338.75  # N97
$ws.Cells.Item(107, 8).Value = 3288.6667  # H107
$ws.Cells.Item(107, 9).Value = 3200  # I107
$ws.Cells.Item(107, 10).Value = 3333  # J107
$ws.Cells.Item(107, 11).Value = 3200  # K107
$ws.Cells.Item(107, 12).Value = 3333  # L107
$ws.Cells.Item(107, 13).Value = -1280  # M107
$ws.Cells.Item(107, 14).Value = -7173  # N107
$ws.Cells.Item(113, 8).Value = 3309.1  # H113
$ws.Cells.Item(113, 9).Value = 2531.8333  # I113
$ws.Cells.Item(113, 11).Value = 2531.8333  # K113
$ws.Cells.Item(113, 13).Value = -361.8332999999998  # M113

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(59, 8).Value = 30000  # H59
$ws.Cells.Item(59, 10).Value = 30000  # J59
$ws.Cells.Item(59, 12).Value = 30000  # L59
$ws.Cells.Item(59, 14).Value = -31308  # N59
$ws.Cells.Item(61, 8).Value = 5575.45  # H61
$ws.Cells.Item(61, 9).Value = 3056  # I61
$ws.Cells.Item(61, 10).Value = 7636.8184  # J61
$ws.Cells.Item(61, 11).Value = 3056  # K61
$ws.Cells.Item(61, 12).Value = 7636.8184  # L61
$ws.Cells.Item(61, 13).Value = -2854  # M61
$ws.Cells.Item(61, 14).Value = -8040.8184  # N61
$ws.Cells.Item(82, 8).Value = 3842.8572  # H82
$ws.Cells.Item(82, 9).Value = 4900  # I82
$ws.Cells.Item(82, 11).Value = 4900  # K82
$ws.Cells.Item(82, 13).Value = -4539  # M82
$ws.Cells.Item(85, 8).Value = 3842.8572  # H85
$ws.Cells.Item(85, 9).Value = 4900  # I85
$ws.Cells.Item(85, 11).Value = 4900  # K85
$ws.Cells.Item(85, 13).Value = -3652  # M85
$ws.Cells.Item(113, 8).Value = 5575.45  # H113
$ws.Cells.Item(113, 9).Value = 3056  # I113
$ws.Cells.Item(113, 10).Value = 7636.8184  # J113
$ws.Cells.Item(113, 11).Value = 3056  # K113
$ws.Cells.Item(113, 12).Value = 7636.8184  # L113
$ws.Cells.Item(113, 13).Value = -886  # M113
$ws.Cells.Item(113, 14).Value = -11976.8184  # N113

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1003.8125  # H100
$ws.Cells.Item(100, 9).Value = 585.8  # I100
$ws.Cells.Item(100, 11).Value = 1171.6  # K100
$ws.Cells.Item(100, 13).Value = -630.5999999999999  # M100
$ws.Cells.Item(132, 8).Value = 3392.9333  # H132
$ws.Cells.Item(132, 9).Value = 2991.5  # I132
$ws.Cells.Item(132, 10).Value = 4998.6665  # J132
$ws.Cells.Item(132, 11).Value = 8974.5  # K132
$ws.Cells.Item(132, 12).Value = 14995.9995  # L132
$ws.Cells.Item(132, 13).Value = -6444.5  # M132
$ws.Cells.Item(132, 14).Value = -20055.9995  # N132
